$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (same
#    layout/fund list) and placing it right before "2022-Q3".
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Update the Q4-specific figures (fund size / stock position / position
# ratio / held market value) - these are stored as text, so force a text
# number format before writing so the stored cell stays a string like the
# rest of the sheet (matches "11.33" style text cells already there).
$q4.Range("D2:G5").NumberFormat = "@"

$q4.Range("D2").Value2 = "11.33"
$q4.Range("E2").Value2 = "90.17"
$q4.Range("F2").Value2 = "6.19"
$q4.Range("G2").Value2 = "0.7013"

$q4.Range("D3").Value2 = "11.33"
$q4.Range("E3").Value2 = "90.17"
$q4.Range("F3").Value2 = "6.19"
$q4.Range("G3").Value2 = "0.7013"

$q4.Range("D4").Value2 = "5.56"
$q4.Range("E4").Value2 = "90.17"
$q4.Range("F4").Value2 = "6.19"
$q4.Range("G4").Value2 = "0.3442"

$q4.Range("D5").Value2 = "-11.33"
$q4.Range("E5").Value2 = "90.17"
$q4.Range("F5").Value2 = "6.19"
$q4.Range("G5").Value2 = "-0.7013"

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a new leading row for
#    2022-Q4 and push the existing quarters down by one row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 8 is brand new - give it the same look (style) as row 7 before
# filling it in.
$summary.Cells.Item(7,1).Copy()
$summary.Cells.Item(8,1).PasteSpecial(-4122)

# Shift existing quarter rows (2..7) down into (3..8), bottom-up so we
# never overwrite data we still need to read.
for ($r = 7; $r -ge 2; $r--) {
    $summary.Cells.Item($r + 1, 2).Value2 = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($r + 1, 3).Value2 = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($r + 1, 4).Value2 = $summary.Cells.Item($r, 4).Value2
}

# Fill in the new 2022-Q4 row at the top of the data.
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 4
$summary.Cells.Item(2, 4).Value2 = 1.05

# Recompute the leading index column (0-based row position).
for ($r = 2; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}

# ------------------------------------------------------------------
# 3) Keep the originally-active tab ("2021-Q2") selected, since it
#    naturally shifted position but shouldn't change which sheet is
#    active.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
